# Refresh the crypto price/volume table to the latest scrape, matching the
# GitHub Actions commit 'Updated cryptos list ... with GitHub Actions'.
# Most rows just get new Price (D) / Volume(1h) (E) text; a few rows were
# re-ordered in the source feed, so Coin (B) and Link (C) are rewritten too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '51.972.35'
$ws.Cells.Item(2, 5).Value = '  +1.56%  '
$ws.Cells.Item(3, 4).Value = '2.819.84'
$ws.Cells.Item(3, 5).Value = '  +1.95%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '352.56'
$ws.Cells.Item(5, 5).Value = '  -0.51%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '113.77'
$ws.Cells.Item(6, 5).Value = '  +5.33%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.562'
$ws.Cells.Item(7, 5).Value = '  +2.41%  '
$ws.Cells.Item(8, 5).Value = '  +0.05%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.622'
$ws.Cells.Item(9, 5).Value = '  +6.79%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '40.54'
$ws.Cells.Item(10, 5).Value = '  +2.86%  '
$ws.Cells.Item(11, 5).Value = '  -0.64%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0847'
$ws.Cells.Item(12, 5).Value = '  +1.62%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '20.00'
$ws.Cells.Item(13, 5).Value = '  -0.17%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '7.83'
$ws.Cells.Item(14, 5).Value = '  +4.05%  '
$ws.Cells.Item(15, 4).Value = '3.260.84'
$ws.Cells.Item(15, 5).Value = '  +1.97%  '
$ws.Cells.Item(16, 2).Value = 'WrappedEther'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(16, 4).Value = '2.832.78'
$ws.Cells.Item(16, 5).Value = '  +2.28%  '
$ws.Cells.Item(17, 2).Value = 'Polygon'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.972'
$ws.Cells.Item(17, 5).Value = '  +4.54%  '
$ws.Cells.Item(18, 4).Value = '51.995.84'
$ws.Cells.Item(18, 5).Value = '  +1.74%  '
$ws.Cells.Item(19, 5).Value = '  +9.75%  '
$ws.Cells.Item(20, 5).Value = '  -1.21%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.64'
$ws.Cells.Item(21, 5).Value = '  +4.01%  '
$ws.Cells.Item(22, 5).Value = '  +1.93%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '70.71'
$ws.Cells.Item(23, 5).Value = '  +1.50%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '269.66'
$ws.Cells.Item(24, 5).Value = '  +1.68%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.77'
$ws.Cells.Item(25, 5).Value = '  +1.92%  '
$ws.Cells.Item(26, 5).Value = '  +1.35%  '
$ws.Cells.Item(27, 5).Value = '  -0.08%  '
$ws.Cells.Item(28, 5).Value = '  +1.07%  '
$ws.Cells.Item(29, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '39.07'
$ws.Cells.Item(29, 5).Value = '  +12.59%  '
$ws.Cells.Item(30, 2).Value = 'Cosmos'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '10.57'
$ws.Cells.Item(30, 5).Value = '  +4.04%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.28'
$ws.Cells.Item(31, 5).Value = '  +1.49%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '52.87'
$ws.Cells.Item(32, 5).Value = '  +2.02%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '6.19'
$ws.Cells.Item(33, 5).Value = '  +2.48%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0907'
$ws.Cells.Item(34, 5).Value = '  +9.32%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.70'
$ws.Cells.Item(35, 5).Value = '  +4.11%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.0454'
$ws.Cells.Item(36, 5).Value = '  +2.32%  '
$ws.Cells.Item(37, 5).Value = '  -0.08%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '19.07'
$ws.Cells.Item(38, 5).Value = '  +4.83%  '
$ws.Cells.Item(39, 5).Value = '  +2.56%  '
$ws.Cells.Item(40, 5).Value = '  +3.57%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.60'
$ws.Cells.Item(41, 5).Value = '  +3.29%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.117'
$ws.Cells.Item(42, 5).Value = '  +2.32%  '
$ws.Cells.Item(43, 2).Value = 'EnergySwap'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '22.46'
$ws.Cells.Item(43, 5).Value = '  +1.66%  '
$ws.Cells.Item(44, 2).Value = 'WEMIXToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.23'
$ws.Cells.Item(44, 5).Value = '  +2.01%  '
$ws.Cells.Item(45, 2).Value = 'Monero'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '121.15'
$ws.Cells.Item(45, 5).Value = '  +0.68%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.52'
$ws.Cells.Item(46, 5).Value = '  +8.80%  '
$ws.Cells.Item(47, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.48'
$ws.Cells.Item(47, 5).Value = '  +9.61%  '
$ws.Cells.Item(48, 2).Value = 'Maker'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(48, 4).Value = '2.140.45'
$ws.Cells.Item(48, 5).Value = '  +2.63%  '
$ws.Cells.Item(49, 5).Value = '  +12.24%  '
$ws.Cells.Item(50, 2).Value = 'BEAM'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0324'
$ws.Cells.Item(50, 5).Value = '  +16.84%  '
$ws.Cells.Item(51, 2).Value = 'TheGraph'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.224'
$ws.Cells.Item(51, 5).Value = '  +18.59%  '
